# Refresh Universalis market-price derived columns (H:N) across the Leve-profit
# sheets (ALC, ARM, BSM, CRP, CUL, WVR) following a scheduled price-data pull.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 2460.0908
$ws.Range("I86").Value = 1507
$ws.Range("J86").Value = 5001.6665
$ws.Range("K86").Value = 1507
$ws.Range("L86").Value = 5001.6665
$ws.Range("M86").Value = -384
$ws.Range("N86").Value = -7247.6665

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 2460.0908
$ws.Range("I89").Value = 1507
$ws.Range("J89").Value = 5001.6665
$ws.Range("K89").Value = 7535
$ws.Range("L89").Value = 25008.3325
$ws.Range("M89").Value = -1919
$ws.Range("N89").Value = -36240.3325

# Row 129: Practical Command
$ws.Range("H129").Value = 1562.6923
$ws.Range("I129").Value = 724.4
$ws.Range("J129").Value = 1762.2858
$ws.Range("K129").Value = 2173.2
$ws.Range("L129").Value = 5286.857400000001
$ws.Range("M129").Value = 2826.8
$ws.Range("N129").Value = -15286.8574

# Row 135: For Tired Minds
$ws.Range("H135").Value = 3963.8445
$ws.Range("I135").Value = 2954.04
$ws.Range("J135").Value = 5226.1
$ws.Range("K135").Value = 26586.36
$ws.Range("L135").Value = 47034.9
$ws.Range("M135").Value = -24051.36
$ws.Range("N135").Value = -52104.9

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2147.7896
$ws.Range("I137").Value = 1344.375
$ws.Range("J137").Value = 2732.0908
$ws.Range("K137").Value = 4033.125
$ws.Range("L137").Value = 8196.2724
$ws.Range("M137").Value = -1483.125
$ws.Range("N137").Value = -13296.2724

# Row 138: All-night Crafting
$ws.Range("H138").Value = 670352.4399999999
$ws.Range("I138").Value = 1535.5555
$ws.Range("J138").Value = 1673577.9
$ws.Range("K138").Value = 4606.666499999999
$ws.Range("L138").Value = 5020733.699999999
$ws.Range("M138").Value = 533.3335000000006
$ws.Range("N138").Value = -5031013.699999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3030.8
$ws.Range("I32").Value = 3017.8362
$ws.Range("K32").Value = 3017.8362
$ws.Range("M32").Value = -2730.8362

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 4512.3076
$ws.Range("I74").Value = 4182.5
$ws.Range("J74").Value = 5040
$ws.Range("K74").Value = 4182.5
$ws.Range("L74").Value = 5040
$ws.Range("M74").Value = -3308.5
$ws.Range("N74").Value = -6788

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 4512.3076
$ws.Range("I77").Value = 4182.5
$ws.Range("J77").Value = 5040
$ws.Range("K77").Value = 20912.5
$ws.Range("L77").Value = 25200
$ws.Range("M77").Value = -16544.5
$ws.Range("N77").Value = -33936

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 11765483
$ws.Range("I102").Value = 907.7857
$ws.Range("J102").Value = 66666830
$ws.Range("K102").Value = 907.7857
$ws.Range("L102").Value = 66666830
$ws.Range("M102").Value = 714.2143
$ws.Range("N102").Value = -66670074

# Row 104: See Shields by the Sea Shore
$ws.Range("H104").Value = 68999.5
$ws.Range("J104").Value = 68999.5
$ws.Range("L104").Value = 68999.5
$ws.Range("N104").Value = -75987.5

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3657.75
$ws.Range("I122").Value = 3657.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10973.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8523.25
$ws.Range("N122").ClearContents()

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3519.5
$ws.Range("I132").Value = 2282.2666
$ws.Range("K132").Value = 6846.7998
$ws.Range("M132").Value = -4316.7998

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3212.743
$ws.Range("I105").Value = 2132.1667
$ws.Range("J105").Value = 4356.8823
$ws.Range("K105").Value = 2132.1667
$ws.Range("L105").Value = 4356.8823
$ws.Range("M105").Value = -385.1667000000002
$ws.Range("N105").Value = -7850.8823

# Row 107: The Gold Experience
$ws.Range("H107").Value = 5459
$ws.Range("I107").Value = 6840
$ws.Range("J107").Value = 2006.5
$ws.Range("K107").Value = 6840
$ws.Range("L107").Value = 2006.5
$ws.Range("M107").Value = -4920
$ws.Range("N107").Value = -5846.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3638.6614
$ws.Range("I31").Value = 1933.6538
$ws.Range("J31").Value = 4870.0557
$ws.Range("K31").Value = 1933.6538
$ws.Range("L31").Value = 4870.0557
$ws.Range("M31").Value = -1638.6538
$ws.Range("N31").Value = -5460.0557

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3638.6614
$ws.Range("I34").Value = 1933.6538
$ws.Range("J34").Value = 4870.0557
$ws.Range("K34").Value = 1933.6538
$ws.Range("L34").Value = 4870.0557
$ws.Range("M34").Value = -1731.6538
$ws.Range("N34").Value = -5274.0557

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 8453.799999999999
$ws.Range("I105").Value = 8623
$ws.Range("K105").Value = 8623
$ws.Range("M105").Value = -6876

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 3396.3333
$ws.Range("I3").Value = 3396.3333
$ws.Range("K3").Value = 10188.9999
$ws.Range("M3").Value = -10076.9999

# Row 56: Culture Club
$ws.Range("H56").Value = 897956.75
$ws.Range("I56").Value = 897956.75
$ws.Range("K56").Value = 897956.75
$ws.Range("M56").Value = -897426.75

# Row 68: Such a Butter Face
$ws.Range("H68").Value = 4022.88
$ws.Range("I68").Value = 1797
$ws.Range("J68").Value = 5070.353
$ws.Range("K68").Value = 5391
$ws.Range("L68").Value = 15211.059
$ws.Range("M68").Value = -4580
$ws.Range("N68").Value = -16833.059

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 4022.88
$ws.Range("I71").Value = 1797
$ws.Range("J71").Value = 5070.353
$ws.Range("K71").Value = 16173
$ws.Range("L71").Value = 45633.177
$ws.Range("M71").Value = -12117
$ws.Range("N71").Value = -53745.177

$ws = $wb.Worksheets.Item("WVR")
# Row 64: Ribbon of Remembrance
$ws.Range("H64").Value = 30001
$ws.Range("I64").Value = 30001
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 30001
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -29753
$ws.Range("N64").ClearContents()

# Row 67: The Road Was a Ribbon of Moonlight (L)
$ws.Range("H67").Value = 30001
$ws.Range("I67").Value = 30001
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 30001
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -29143
$ws.Range("N67").ClearContents()

# Row 107: Flax Wax
$ws.Range("H107").Value = 1085.3334
$ws.Range("I107").Value = 999
$ws.Range("J107").Value = 1102.6
$ws.Range("K107").Value = 2997
$ws.Range("L107").Value = 3307.8
$ws.Range("M107").Value = -1077
$ws.Range("N107").Value = -7147.799999999999

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 4509.3057
$ws.Range("I136").Value = 1985
$ws.Range("J136").Value = 5620
$ws.Range("K136").Value = 5955
$ws.Range("L136").Value = 16860
$ws.Range("M136").Value = -3405
$ws.Range("N136").Value = -21960
